$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Shape 2 (id=3, "Content Placeholder 2") holds the bullet list being rewritten.
$sh = $s.Shapes.Item(2)

# Resize / reposition per the new layout (wider box spanning more of the slide).
$sh.Left = 838199 / 12700
$sh.Top = 1825624 / 12700
$sh.Width = 10515600 / 12700
$sh.Height = 4667251 / 12700

$tf = $sh.TextFrame
$tr = $tf.TextRange

# Replace the old "Cons"-style critique bullets with the new "Pros" bullets.
$tr.Text = "Paper was well structured and easy to follow`r" + `
           "Objective tree was well thought out and had good visuals`r" + `
           "Technical feasibility was informative on the tasks at hand for the team`r" + `
           "Time feasibility set short term goals as well, long-term goals`r" + `
           "Requirement Specifications was thorough`r" + `
           "This section really emphasized more goals for the team while also introducing the board requirements`r" + `
           "This section had some points that were well stated, an example would be identifying the blocks`r"

# Paragraph 1: "Paper was well structured and easy to follow" - top level, 24pt
$para = $tr.Paragraphs(1,1)
$para.Font.Size = 24

# Paragraph 2: "Objective tree was well thought out and had good visuals" - top level, 24pt
$para = $tr.Paragraphs(2,1)
$para.Font.Size = 24

# Paragraph 3: "Technical feasibility was informative on the tasks at hand for the team" - top level, 24pt
$para = $tr.Paragraphs(3,1)
$para.Font.Size = 24

# Paragraph 4: "Time feasibility set short term goals as well, long-term goals" - sub level, 20pt
$para = $tr.Paragraphs(4,1)
$para.IndentLevel = 2
$para.Font.Size = 20

# Paragraph 5: "Requirement Specifications was thorough" - top level, 24pt
$para = $tr.Paragraphs(5,1)
$para.Font.Size = 24

# Paragraph 6: "This section really emphasized..." - sub level, 20pt
$para = $tr.Paragraphs(6,1)
$para.IndentLevel = 2
$para.Font.Size = 20

# Paragraph 7: "This section had some points..." - sub level, 20pt
$para = $tr.Paragraphs(7,1)
$para.IndentLevel = 2
$para.Font.Size = 20

# Trailing empty paragraph (endParaRPr sz=2400) after the last bullet.
$para = $tr.Paragraphs(8,1)
$para.Font.Size = 24

# Remove the second content box (id=5, "Alternative and Tradeoffs" / "Feasibility Assessment").
$s.Shapes.Item(3).Delete()
